$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before row 243, shifting the existing rows
# 243:266 down to 246:269. This mirrors the diff, where three brand-new
# weekly records (Clementina, 31-05-2023, Region de O'Higgins) were
# inserted at the top of this date-ordered block and all subsequent
# rows moved down by three positions.
$ws.Range("A243:A245").EntireRow.Insert()

# New row 243: Clementina / Especial
$ws.Cells.Item(243, 1).Value2 = 7
$ws.Cells.Item(243, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(243, 3).Value2 = "Ñuble"
$ws.Cells.Item(243, 4).Value2 = 45077
$ws.Cells.Item(243, 5).Value2 = 16
$ws.Cells.Item(243, 6).Value2 = "Fruta"
$ws.Cells.Item(243, 7).Value2 = 100102
$ws.Cells.Item(243, 8).Value2 = "Cítricos"
$ws.Cells.Item(243, 9).Value2 = 100102004
$ws.Cells.Item(243, 10).Value2 = "Mandarina"
$ws.Cells.Item(243, 11).Value2 = "Clementina"
$ws.Cells.Item(243, 12).Value2 = "Especial"
$ws.Cells.Item(243, 13).Value2 = 50
$ws.Cells.Item(243, 14).Value2 = 13000
$ws.Cells.Item(243, 15).Value2 = 13000
$ws.Cells.Item(243, 16).Value2 = 13000
$ws.Cells.Item(243, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(243, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(243, 19).Value2 = 1300
$ws.Cells.Item(243, 20).Value2 = 10

# New row 244: Clementina / Primera
$ws.Cells.Item(244, 1).Value2 = 7
$ws.Cells.Item(244, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(244, 3).Value2 = "Ñuble"
$ws.Cells.Item(244, 4).Value2 = 45077
$ws.Cells.Item(244, 5).Value2 = 16
$ws.Cells.Item(244, 6).Value2 = "Fruta"
$ws.Cells.Item(244, 7).Value2 = 100102
$ws.Cells.Item(244, 8).Value2 = "Cítricos"
$ws.Cells.Item(244, 9).Value2 = 100102004
$ws.Cells.Item(244, 10).Value2 = "Mandarina"
$ws.Cells.Item(244, 11).Value2 = "Clementina"
$ws.Cells.Item(244, 12).Value2 = "Primera"
$ws.Cells.Item(244, 13).Value2 = 40
$ws.Cells.Item(244, 14).Value2 = 12000
$ws.Cells.Item(244, 15).Value2 = 12000
$ws.Cells.Item(244, 16).Value2 = 12000
$ws.Cells.Item(244, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(244, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(244, 19).Value2 = 1200
$ws.Cells.Item(244, 20).Value2 = 10

# New row 245: Clementina / Segunda
$ws.Cells.Item(245, 1).Value2 = 7
$ws.Cells.Item(245, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(245, 3).Value2 = "Ñuble"
$ws.Cells.Item(245, 4).Value2 = 45077
$ws.Cells.Item(245, 5).Value2 = 16
$ws.Cells.Item(245, 6).Value2 = "Fruta"
$ws.Cells.Item(245, 7).Value2 = 100102
$ws.Cells.Item(245, 8).Value2 = "Cítricos"
$ws.Cells.Item(245, 9).Value2 = 100102004
$ws.Cells.Item(245, 10).Value2 = "Mandarina"
$ws.Cells.Item(245, 11).Value2 = "Clementina"
$ws.Cells.Item(245, 12).Value2 = "Segunda"
$ws.Cells.Item(245, 13).Value2 = 30
$ws.Cells.Item(245, 14).Value2 = 10000
$ws.Cells.Item(245, 15).Value2 = 10000
$ws.Cells.Item(245, 16).Value2 = 10000
$ws.Cells.Item(245, 17).Value2 = "$/bandeja 10 kilos"
$ws.Cells.Item(245, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(245, 19).Value2 = 1000
$ws.Cells.Item(245, 20).Value2 = 10
